# Commit: "Sat, Jun 20, 2020  2:04:50 PM"
#
# The canonical-XML diff shows the whole <a:clrScheme> (and the a:theme /
# a:clrScheme @name) of ppt/theme/theme1.xml (the theme used by the slide
# master, i.e. the presentation's visible "Design") and
# ppt/theme/theme2.xml (the theme used only by the notes master) being
# swapped: theme1 flips from the "Integral" palette to the stock "Office"
# palette, and theme2 flips from "Office" to "Integral". Everything else
# (fontScheme/fmtScheme) is already byte-identical between the two parts,
# so the only real effect of the edit is a 12-colour palette change on the
# presentation's main theme (what the Design/Colours gallery shows).
#
# PowerPoint's automation model exposes that palette as
# Master.Theme.ThemeColorScheme (12 items, 1-based, in the fixed order
# dk1, lt1, dk2, lt2, accent1-6, hlink, folHlink) with a settable .RGB on
# each item - so we repaint the master's colour scheme to the "Office"
# values from the diff.

$p = $ppt.ActivePresentation
$master = $p.SlideMaster
$scheme = $master.Theme.ThemeColorScheme

# Target palette ("Office Theme" clrScheme from the diff), in VBA RGB()
# (R + G*256 + B*65536) form, ordered dk1,lt1,dk2,lt2,accent1-6,hlink,folHlink
$officeColors = @(
    0,          # dk1      000000
    16777215,   # lt1      FFFFFF
    6968388,    # dk2      44546A
    15132391,   # lt2      E7E6E6
    13998939,   # accent1  5B9BD5
    3243501,    # accent2  ED7D31
    10855845,   # accent3  A5A5A5
    49407,      # accent4  FFC000
    12874308,   # accent5  4472C4
    4697456,    # accent6  70AD47
    12673797,   # hlink    0563C1
    7491477     # folHlink 954F72
)

for ($i = 1; $i -le 12; $i++) {
    $scheme.Item($i).RGB = $officeColors[$i - 1]
}
